$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (District). This shifts the existing
# column F (District) and everything after it one column to the right,
# turning old column F into column G, and leaves a blank column F ready
# to hold the new "Address" data.
$ws.Columns("F").Insert()

# Header labels for the new Address column
$ws.Range("F2").Value = "Address"
$ws.Range("F3").Value = ""

# Per-row Address values (school name / location), mirrors the second line
# of column B minus the trailing district name.
$addresses = @{
    4 = "Govt. High School Thirthahalli"
    5 = "G H S NitturHosanagar"
    6 = "J P N High School"
    7 = "Paper Town High School Paper Town Bhadravathi"
    8 = "G H S Kadasuru Sorab"
    9 = "G H S BilikiShikaripura"
    10 = "G H S Thammadihalli"
    11 = "G H S Kanale Sagar"
    12 = "Sri Parameshwara High School Muduba SiddapuraShikaripura"
    13 = "G H S Kattinakaru Sagara"
    14 = "G U H S Soraba"
    15 = "Kasturba Girls High School"
    16 = "G H S HunshanakatteThirthahalli"
    17 = "G H S Kaginalli Shikaripura"
    18 = "G H S HunasekatteBhadravathi"
    19 = "Bapuji High School Shikaripura"
    20 = "G J C Shikaripura"
    21 = "G H S Kuskur(RMSA)Shikaripura"
    22 = "Govt. Girls High School AnavattiSorab"
    23 = "G H S BalemaranahalliBhadravathi"
    24 = "Sri Channabasaveshwara High School"
    25 = "Kasturba Girls High School"
    26 = "G J C Amrutha Hosanagara"
    27 = "G H S Sirigere"
    28 = "G J C MalalimakkiThirthahalli"
    29 = "G J C HolehonnurBhadravathi"
    30 = "G G P U C New town Bhadravathi"
    31 = "G H S Tadagalale Sagar"
    32 = "Govt. High School Durgigudi"
    33 = "G H S ThanikalThirthahalli"
    34 = "G H S Subhash NagarSagar"
    35 = "G H S Kommanalu"
    36 = "Dr. U R AnanthamurthyG H S Thirthahalli"
    37 = "G H S BavikaisaruThirthahalli"
    38 = "National P U CollegeB R ProjectBhadravathi"
    39 = "V S High School NisaraniSorab"
    40 = "G H S ArasaluHosanagar"
    41 = "Channammaji High School AlagerimandriHosanagara"
    42 = "G H S Kadekal"
    43 = "G H S Halesoraba Sorab"
    44 = "G H S BarurSagar"
    45 = "G P U C High School SectionMegaravalliThirthahalli"
    46 = "Govt. High School SonaleHosanagara"
    47 = "G H S Mavinakere Colony Bhadravathi"
    48 = "G H S Hosur GuddekeriThirthahalli"
    49 = "G H S Basavani Thirthahalli"
    50 = "G H S Mandaghatta"
    51 = "Kalleshwara High School YalageriShikaripur"
    52 = "G H S Yadur Hosanagara"
    53 = "G P U T K RoadBhadravathi"
    54 = "G H S Tumari"
    55 = "G H S KalmaneShikaripura"
    56 = "G H S Nagara Hosanagara"
}

foreach ($row in $addresses.Keys) {
    $ws.Range("F$row").Value = $addresses[$row]
}
